# "more tweak after presentation"
#  1. Add alt-text ("descr") to the two pictures on the last existing slide.
#  2. Append a new slide containing a small "Copyright (c)" textbox.

$p = $ppt.ActivePresentation

# --- 1. Alt text for the two pictures on slide 9 --------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(3).AlternativeText = "A screenshot of a cell phone`n`nDescription automatically generated"
$s9.Shapes.Item(5).AlternativeText = "A purple flower in a field`n`nDescription automatically generated"

# --- 2. New trailing slide with a copyright textbox ------------------------
# Duplicate the last slide (rather than Slides.Add) so the new slide picks
# up the normal per-slide boilerplate (group transform, color-map override,
# etc.) that a brand new blank slide wouldn't otherwise carry, then strip
# its copied shapes/background back out before adding the real content.
$null = $s9.Duplicate()
$s = $p.Slides.Item($p.Slides.Count)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $s.Shapes.Item($i).Delete()
}
$s.FollowMasterBackground = $true

$left = 9164055 / 12700.0
$top = 149088 / 12700.0
$width = 3027945 / 12700.0
$height = 369332 / 12700.0

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 1"
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $false
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "Copyright " + [char]0x00A9 + " 2019 Xiangshi Yin"
